$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (row 1) to new machine-readable names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case Spanish connector words (de, del, la, las, los, el, y, en)
#    inside the state/municipality name columns (A and B), rows 2-778.
for ($r = 2; $r -le 778; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value()
    if ($valA -ne $null -and $valA -is [string]) {
        $newA = $valA -replace '\bde\b','De' -replace '\bdel\b','Del' -replace '\blas\b','Las' -replace '\bla\b','La' -replace '\blos\b','Los' -replace '\bel\b','El' -replace '\by\b','Y' -replace '\ben\b','En'
        $cellA.Value = $newA
    }
    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value()
    if ($valB -ne $null -and $valB -is [string]) {
        $newB = $valB -replace '\bde\b','De' -replace '\bdel\b','Del' -replace '\blas\b','Las' -replace '\bla\b','La' -replace '\blos\b','Los' -replace '\bel\b','El' -replace '\by\b','Y' -replace '\ben\b','En'
        $cellB.Value = $newB
    }
}

# 3. Fix floating point precision difference in D332 (Tamazula De Gordiano)
$ws.Range("D332").Value = 0.009484536082474229

# 4. Delete footer/metadata rows 780-784 (row 779 was already a gap row)
$ws.Rows("780:784").Delete()
